$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.884.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.841.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.69%  "
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4669"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3660"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07132"
$ws.Range("D9").Style = "Normal"
$ws.Range("E10").Value = "  +2.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07699"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.859.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.279"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.388"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.54%  "
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008619"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.901.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.015"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.929"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.027"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.884"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08861"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.207"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7459"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("E33").Value = "  +5.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.780"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.082"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01942"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.969"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5197"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.899"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1514"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.83%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.120"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4691"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06026"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8861"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.05%  "
